$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the power value for row 2 (ID_HeatingElement = 1)
$ws.Range("B2").Value = 100000

# Remove row 3 (ID_HeatingElement = 2) entirely - only one heating element remains
$ws.Rows.Item(3).Delete()

# Update the active selection to reflect the new extent of data
[void]$ws.Range("B3").Select()
